$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $ok = $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Host "WARNING: replace failed for: $old"
    }
}

function Delete-Text($old) {
    $ok = $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, "", 2)
    if (-not $ok) {
        Write-Host "WARNING: delete failed for: $old"
    }
}

# ------------------------------------------------------------------
# Title
# ------------------------------------------------------------------
Replace-Text "The Timeless Melody of Music" "Exploring the Marvels of the Human Body: A Journey through Biology"

# ------------------------------------------------------------------
# Author name: "Isabella Sinclair" -> "Dr. Ella McPherson"
# ------------------------------------------------------------------
Replace-Text "Isabella Sinclair" "Dr. Ella McPherson"

# ------------------------------------------------------------------
# Email: "isabella.sinclair@melodicinstitute.edu" -> "emcpherson@bioacademy.edu"
# ------------------------------------------------------------------
Replace-Text "isabella.sinclair@melodicinstitute" "emcpherson@bioacademy"

# ------------------------------------------------------------------
# Intro paragraph (first body paragraph, sz 24)
# ------------------------------------------------------------------
Replace-Text "Music, a universal language that transcends boundaries and cultures, has captivated humanity for centuries" "Every breath we take, every beat of our heart, every thought that crosses our minds - all are intricate manifestations of the wonders of biology"

Replace-Text " From the haunting melodies of ancient civilizations to the modern-day symphonies, music's power to evoke emotions, communicate stories, and connect people is undeniable" " Biology, a science that delves into the world of living organisms, unveils the captivating tapestry of life's processes"

Replace-Text " In the realm of sound, music weaves a tapestry of experiences, shaping our understanding of the world around us" " In the symphony of existence, biology plays a key role, conducting the harmonious dance of molecules and cells, tissues and organs"

Replace-Text " Its influence extends far beyond mere entertainment; it serves as a catalyst for social change, a healer of wounds, and a source of profound inspiration" " As high school students embark on this enthralling journey, they will become explorers of the enigmatic realm of life, unraveling the intricate web of biological phenomena"

Replace-Text "Music has long been intertwined with human history" "In this exploration, we will delve into the depths of our bodies, delving into the intricacies of our cells, the building blocks of life"

Replace-Text " Cave paintings, ancient instruments, and oral traditions reveal the deep-rooted connection between our species and the art of sound" " From exploring the molecular dance of DNA to deciphering the complex communications between cells, our journey will unravel the mysteries of human physiology"

Replace-Text " From tribal rituals to elaborate courtly performances, music has played an integral role in shaping cultures, forging identities, and preserving traditions" " With each step, students will discover the elegance of adaptation, the resilience of life, and the interconnectedness of all living things"

# Runs fully removed: ". It provides a means of self-expression, enabling individuals to
# communicate their innermost thoughts and emotions without words. Moreover, music fosters
# a sense of community, uniting people from diverse backgrounds in shared experiences of
# joy, sorrow, and reflection."
Delete-Text ". It provides a means of self-expression, enabling individuals to communicate their innermost thoughts and emotions without words. Moreover, music fosters a sense of community, uniting people from diverse backgrounds in shared experiences of joy, sorrow, and reflection."

Replace-Text "In modern times, music has evolved into a multifaceted phenomenon" "Biology extends far beyond the realm of human existence, encompassing the diversity and abundance of life on Earth"

Replace-Text " The advent of recording technology and the rise of mass media have transformed the way music is created, distributed, and consumed" " From the depths of the oceans to the soaring heights of rainforests, students will decipher the symbiotic relationships that shape ecosystems, marveling at the intricate dance of predator and prey"

# Runs fully removed: ". Today, we have access to an unprecedented diversity of musical
# genres, styles, and artists, accessible at our fingertips through streaming services and
# digital platforms."
Delete-Text ". Today, we have access to an unprecedented diversity of musical genres, styles, and artists, accessible at our fingertips through streaming services and digital platforms."

Replace-Text " The internet has also facilitated the emergence of virtual communities where music lovers can connect, share their passion, and collaborate on creative projects" " By unlocking the secrets of evolution, they will gain a profound appreciation for the unity and diversity of life on our planet, fostering a sense of wonder and responsibility for the natural world that sustains us"

# ------------------------------------------------------------------
# Summary paragraph
# ------------------------------------------------------------------
Replace-Text "Music's impact on humanity is profound and everlasting" "Biology is an awe-inspiring voyage through the marvels of life, revealing the extraordinary complexities of living organisms"

Replace-Text " It transcends time and cultures, serving as a universal language that speaks to the human soul" " Embracing a holistic approach, this journey delves into the intricate workings of human physiology, uncovers the principles of evolution and adaptation, and celebrates the breathtaking diversity of ecosystems"

# Runs fully removed: " From ancient civilizations to modern-day societies, music has played
# a vital role in shaping history, culture, and human connection."
Delete-Text " From ancient civilizations to modern-day societies, music has played a vital role in shaping history, culture, and human connection."

Replace-Text " Whether it be through its ability to evoke emotions, communicate stories, or " " This exploration "

Replace-Text "inspire change, music remains an essential part of our lives" "cultivates scientific curiosity, nurtures an appreciation for the interconnectedness of life, and inspires a commitment to preserving the natural world"

# Runs fully removed: ". It enriches our experiences, heals our wounds, and connects us to
# one another in ways that words cannot."
Delete-Text ". It enriches our experiences, heals our wounds, and connects us to one another in ways that words cannot."

# ------------------------------------------------------------------
# Append a new empty paragraph at the end of the body (before sectPr)
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$endRange = $lastPara.Range
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
